$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- A2 / A3 hold the run date as TEXT (not an Excel date serial). Typing
# "2025-12-06" straight into .Value gets auto-recognized as a date, which
# would stamp a date number-format onto the cell. Route it through a
# formula-result + paste-values round trip instead so the cell keeps its
# original (default/no) style while still landing as literal text.
$helper = $ws.Range("Z1")
$helper.Formula = '="2025-12-06"'
$helper.Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues
$helper.Copy()
$ws.Range("A3").PasteSpecial(-4163)  # xlPasteValues
$helper.ClearContents()
$excel.CutCopyMode = $false

# Row 2: Oklo Inc. (OKLO)
$ws.Range("D2").Value = 105.45
$ws.Range("E2").Value = 55.6
$ws.Range("F2").Value = 15.4
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 76
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 63.9
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 51.54219175917372
$ws.Range("O2").Value = "⚪ 중립 구간"

# Row 3: NuScale Power Corporation (SMR)
$ws.Range("D3").Value = 22.15
$ws.Range("F3").Value = 10.75
$ws.Range("H3").Value = 56
$ws.Range("I3").Value = 66
$ws.Range("J3").Value = 56
$ws.Range("K3").Value = 56.9
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 51.54219175917372
$ws.Range("O3").Value = "⚪ 중립 구간"
